# Commit message: "fixed error: converted m to ft"
# The "results" sheet had leftover name labels in column B (copied over from
# the "mu" sheet) that no longer make sense once the units were converted
# from meters to feet. Clear the stray text values from B4:B11 on the
# "results" sheet, leaving their formatting intact, and update the active
# selection on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("results")

# Clear the (now-incorrect) text values in B4:B11 while preserving styles/formatting.
$ws.Range("B4:B11").ClearContents()

# Activate the results sheet and set the new selection (C3) to match the diff.
$ws.Activate()
$ws.Range("C3").Select()
